# Daily auto push: insert a new log row for 2026/01/18 (日) that got
# appended/merged into the sheet after it had already rolled over into the
# 2026/12/29 block. This shifts every existing row from 662 downward by
# one, and adds the new row 662 with the "13" reading that used to open
# the 2026/12/29 (火) sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 662..703 down to 663..704, leaving a blank row 662 behind.
$ws.Rows("662:662").Insert()

# Column A holds a date-like string ("2026/01/18") that must stay plain
# text (as every other row in the sheet does) instead of being
# auto-coerced into a real Excel date serial by the COM value setter.
$ws.Cells.Item(662, 1).NumberFormat = "@"
$ws.Cells.Item(662, 1).Value = "2026/01/18"
$ws.Cells.Item(662, 2).Value = "日"
$ws.Cells.Item(662, 3).Value = 13
$ws.Cells.Item(662, 4).Value = 201
